$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts old row2..row104 down to row3..row105)
$ws.Rows.Item(2).Insert()

# New row 2 mirrors the values that are now in row 3 (i.e. the old row 2),
# except for the Date column which advances to the new day.
$ws.Cells.Item(2, 1).Value = "18-11-2025"
$ws.Cells.Item(2, 2).Value = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(3, 3).Value2
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(3, 4).Value2
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(3, 5).Value2
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(3, 6).Value2
